$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "http://ontology.deic.dk/cv/M4M-MIK-attempt1/"
$ws.Range("C3").Value = "http://purl.org/M4M-MIK-attempt1/variables/"
$ws.Range("B10").Value = "Agro Ontology"
$ws.Range("B11").Value = "Potato Blight disease observation"
$ws.Range("B12").Value = "https://orcid.org/0000-0002-6384-8944"
$ws.Range("C12").Value = "Minka Karaivanova"
$ws.Range("A13").Value = "dct:creator"
$ws.Range("B13").Value = "https://orcid.org/0000-0003-4093-2147"
$ws.Range("C13").Value = "Eva Overby Bach"
$ws.Range("A14").Value = "dct:creator"
$ws.Range("B14").Value = "https://orcid.org/0000-0002-0721-551X"
$ws.Range("C14").Value = "Ying Wang"
$ws.Range("A15").Value = "dct:creator"
$ws.Range("C15").Value = ""
$ws.Range("A16").Value = "dct:rights"
$ws.Range("B16").Value = "https://spdx.org/licenses/CC-BY-4.0.html"
$ws.Range("C16").Value = "License under which the vocabulary is provided"
$ws.Range("A17").Value = "pav:version"
$ws.Range("B17").Value = "0.0.1"
$ws.Range("C17").Value = "Vocabulary version"
$ws.Range("A18").Value = "pav:createdOn"
$ws.Range("B18").Value = "2022-06-01T19:03:28Z"
$ws.Range("C18").Value = "Date when vocabulary was initially created (follow https://en.wikipedia.org/wiki/ISO_8601)"
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""
$ws.Range("J18").Value = ""
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("A19").Value = "pav:lastUpdatedOn"
$ws.Range("B19").Value = "2022-06-02T19:03:28Z"
$ws.Range("C19").Value = "Date of the last vocabulary update"
$ws.Range("A20").Value = "Definition of terms (optionally properties)"
$ws.Range("B20").Value = ""
$ws.Range("A21").Value = "Identifier"
$ws.Range("B21").Value = "skos:prefLabel@en"
$ws.Range("C21").Value = "qudt:unit(separator=`",`")"
$ws.Range("D21").Value = "skos:altLabel(separator=`";`")"
$ws.Range("E21").Value = "skos:definition@en"
$ws.Range("F21").Value = "dct:source(separator=`",`")"
$ws.Range("G21").Value = "skos:broader(separator=`",`")"
$ws.Range("H21").Value = "skos:exactMatch(separator=`",`")"
$ws.Range("I21").Value = "skos:closeMatch(separator=`",`")"
$ws.Range("J21").Value = "skos:editorialNote@en"
$ws.Range("K21").Value = "dct:creator(separator=`",`")"
$ws.Range("L21").Value = "dct:contributor(separator=`",`")"
$ws.Range("A22").Value = "vars:SampleID"
$ws.Range("B22").Value = "SampleID"
$ws.Range("E22").Value = "Number given to the sample after institutes nameing standard"
$ws.Range("A23").Value = "vars:ObservationID"
$ws.Range("B23").Value = "ObservationID"
$ws.Range("E23").Value = "Running number"
$ws.Range("A24").Value = "vars:CropSeasonYear"
$ws.Range("B24").Value = "CropSeasonYear"
$ws.Range("E24").Value = "4 digit number representing the year the disease was observed. Automatically created from the Blight Tracker App."
$ws.Range("A25").Value = "vars:CountryCode"
$ws.Range("B25").Value = "CountryCode"
$ws.Range("E25").Value = "2 char code as defined in the ISO standard ISO 3166-1 alpha-2 codes "
$ws.Range("F25").Value = "https://www.iso.org/obp/ui/#search`n"
$ws.Range("A26").Value = "vars:GrowthStageName"
$ws.Range("B26").Value = "GrowthStageName"
$ws.Range("F26").Value = "https://en.wikipedia.org/wiki/BBCH-scale_(potato)"
$ws.Range("A27").Value = "vars:SeverityCategoryName"
$ws.Range("B27").Value = "SeverityCategoryName"
$ws.Range("A28").Value = "vars:Disease"
$ws.Range("B28").Value = "Disease"
$ws.Range("E28").Value = "Crop Disease name"
$ws.Range("A29").Value = "vars:PotatoLateBlight"
$ws.Range("B29").Value = "PotatoLateBlight"
$ws.Range("E29").Value = "Potato disease caused by the oomycete Phytophthora infestans"
$ws.Range("G29").Value = "vars:Disease"
$ws.Range("A30").Value = "vars:PotatoEarlyBlight"
$ws.Range("B30").Value = "PotatoEarlyBlight"
$ws.Range("E30").Value = "Potato disease caused by the fungus Alternaria solani"
$ws.Range("G30").Value = "vars:Disease"
$ws.Range("A90").Value = "vars:"
$ws.Range("A91").Value = "vars:"
$ws.Range("A92").Value = "vars:"
